# LOQ4064.xlsx content fix-up
# - fills in the Objetivos / Docentes responsaveis / Programa resumido /
#   Programa / Metodo / Criterio / Norma de recuperacao / Bibliografia
#   cells that were previously blank or carried a leftover value from a
#   neighbouring field, inserting a new row for "Docentes responsaveis"
#   along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Objetivos: (row 10) - replace the stray professor name with the
#    actual Portuguese objectives text.
# ---------------------------------------------------------------------
$objetivosPt = @"
1 - Consolidação e aplicação dos conhecimentos adquiridos em cada uma das áreas específicas do curso de Engenharia Química.  2  Integração dos conhecimentos da Engenharia Química.
"@
$objetivosPt = $objetivosPt.TrimEnd("`r","`n")
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# ---------------------------------------------------------------------
# 2. Insert a new row at 13 for "Docentes responsaveis:" data (shifts the
#    old rows 13-24 down to 14-25).
# ---------------------------------------------------------------------
$ws.Rows("13").Insert()

# The insert leaves a stray, value-less A13 cell carrying column A's bold
# style; strip that cell back to "no style" (copy formatting from a
# genuinely untouched cell) so row 13 ends up with only B13/C13.
$ws.Range("Z1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").ClearContents()

# Give B13/C13 the normal body-text formatting (copy from column B/C of an
# existing data row) and fill in the responsible professor.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$docente = "5816812 - João Paulo Alves Silva"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# ---------------------------------------------------------------------
# 3. Programa resumido: (now row 14) - replace "Semestral" with the short
#    syllabus text in Portuguese.
# ---------------------------------------------------------------------
$resumidoPt = @"
Diagramas para estudos de processos químicos. 2  Estrutura e síntese de processos químicos industriais. 3  Análise de desempenho de processos químicos. 4  Estudo de planta química industrial
"@
$resumidoPt = $resumidoPt.TrimEnd("`r","`n")
$ws.Range("B14").Value = $resumidoPt
$ws.Range("C14").Value = $resumidoPt

# ---------------------------------------------------------------------
# 4. Programa: (now row 16) - fill in the full Portuguese syllabus text
#    (previously held a leftover date value).
# ---------------------------------------------------------------------
$programaPt = @"
1 - Diagramas para estudos de processos químicos: diagramas de bloco; Fluxogramas de processo (PFD); Fluxogramas de instrumentação e tubulação (P&ID).
2  Estrutura e síntese de processos químicos industriais: Hierarquia no planejamento de processos; Etapa 1- Descontínuo ou contínuo; Etapa 2 - Estrutura de entrada/saída de processo; Etapa 3- Estrutura de reciclo; 
3  Análise de desempenho de processos químicos: Modelo de entrada e saída; Ferramentas para a avaliação de processos.
4  Estudo de planta química industrial.
"@
$programaPt = $programaPt.TrimEnd("`r","`n")
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# ---------------------------------------------------------------------
# 5. Método: (now row 19) - the evaluation method text.
# ---------------------------------------------------------------------
$metodo = "Provas escritas e Apresentação de Trabalhos"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# ---------------------------------------------------------------------
# 6. Critério: (now row 20) - the grading criterion text.
# ---------------------------------------------------------------------
$criterio = @"
A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.
"@
$criterio = $criterio.TrimEnd("`r","`n")
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# ---------------------------------------------------------------------
# 7. Norma de recuperação: (now row 21) - make-up exam rule.
# ---------------------------------------------------------------------
$norma = "Média Final = (N + Prova Recuperação)/2"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# ---------------------------------------------------------------------
# 8. Bibliografia: (now row 22) - reference list (previously missing).
# ---------------------------------------------------------------------
$biblio = @"
PERLINGEIRO, Carlos A. G. Engenharia de processos: análise, simulação, otimização e síntese de processos químicos.  Editora Blucher, 2005.
TURTON, BAILIE; WHITING; SHAEIWITZ  Analysis, Synthesis, and Design of Chemical Processes. 3. Ed. LTC Editora, 2005.
COULSON, J. M.; RICHARDSON, J.F. Chemical Engineering Design: Chemical Engineering Volume 6. Editora Fourth, 2005.
HIMMELBLAU, David M. Engenharia química princípios e cálculos. LTC Editora, 2006.
FELDER, R.M; Rousseau, R.W. Princípios elementares dos processos químicos. LTC Editora, 2005.
"@
$biblio = $biblio.TrimEnd("`r","`n")
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
